$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D (Fecha), I (Calidad), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado),
# P (Precio $/Kg). Only rows whose values actually change are listed;
# this reflects a reshuffle of records across dates in the source data.

$rows = @(
    @{ Row = 2;  D = 44799; I = "Primera"; J = 800;  K = 1000; L = 1200; M = 1100; P = 1100 },
    @{ Row = 3;  D = 44278; I = "Segunda"; J = 700;  K = 600;  L = 700;  M = 650;  P = 650 },
    @{ Row = 4;  D = 44278; I = "Tercera"; J = 400;  K = 500;  L = 600;  M = 550;  P = 550 },
    @{ Row = 5;  D = 44201; I = "Segunda"; J = 500;  K = 800;  L = 900;  M = 850;  P = 850 },
    @{ Row = 8;  D = 44544; I = "Primera"; J = 1000; K = 600;  L = 650;  M = 625;  P = 625 },
    @{ Row = 11; D = 44658; I = "Segunda"; J = 1000; K = 600;  L = 650;  M = 625;  P = 625 },
    @{ Row = 12; D = 44229; I = "Segunda"; J = 760;  K = 550;  L = 600;  M = 575;  P = 575 },
    @{ Row = 13; D = 44935; I = "Segunda"; J = 1000; K = 400;  L = 500;  M = 460;  P = 460 },
    @{ Row = 14; D = 44210; I = "Segunda"; J = 900;  K = 600;  L = 700;  M = 650;  P = 650 },
    @{ Row = 15; D = 44573; I = "Tercera"; J = 800;  K = 600;  L = 650;  M = 625;  P = 625 },
    @{ Row = 16; D = 44245; I = "Primera"; J = 800;  K = 850;  L = 900;  M = 875;  P = 875 },
    @{ Row = 17; D = 44245; I = "Segunda"; J = 1000; K = 750;  L = 800;  M = 775;  P = 775 },
    @{ Row = 18; D = 44874; I = "Tercera"; J = 1200; K = 450;  L = 500;  M = 475;  P = 475 },
    @{ Row = 19; D = 44224; I = "Segunda"; J = 800;  K = 850;  L = 900;  M = 875;  P = 875 },
    @{ Row = 20; D = 44474; I = "Segunda"; J = 200;  K = 600;  L = 700;  M = 650;  P = 650 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value  = $r.D   # D: Fecha
    $ws.Cells.Item($row, 9).Value  = $r.I   # I: Calidad
    $ws.Cells.Item($row, 10).Value = $r.J   # J: Volumen
    $ws.Cells.Item($row, 11).Value = $r.K   # K: Precio minimo
    $ws.Cells.Item($row, 12).Value = $r.L   # L: Precio maximo
    $ws.Cells.Item($row, 13).Value = $r.M   # M: Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $r.P   # P: Precio $/Kg
}
